# Update Vip-Vipr2 NATMI output: add new "Inflammatory-Mac" target cluster row
# and refresh all TPM-derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: target cluster = ECs -------------------------------------------------
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7349876666666667
$ws.Range("H2").Value = 2.204963
$ws.Range("M2").Value = 0.077915
$ws.Range("N2").Value = 0.233745
$ws.Range("O2").Value = 0.01344156357222124
$ws.Range("P2").Value = 0.01898188492565723
$ws.Range("Q2").Value = 0.05726656404833334
$ws.Range("R2").Value = 0.515399076435
$ws.Range("S2").Value = 0.01344156357222124
$ws.Range("T2").Value = 0.01898188492565723

# --- Row 3: target cluster = FAPs -------------------------------------------------
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7349876666666667
$ws.Range("H3").Value = 2.204963
$ws.Range("O3").Value = 0.1057841766036689
$ws.Range("P3").Value = 0.1493861228611831
$ws.Range("Q3").Value = 0.4506839023767779
$ws.Range("R3").Value = 4.056155121391001
$ws.Range("S3").Value = 0.1057841766036689
$ws.Range("T3").Value = 0.1493861228611831

# --- Row 4: target cluster changes from MuSCs to new "Inflammatory-Mac" ----------
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7349876666666667
$ws.Range("H4").Value = 2.204963
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.029863
$ws.Range("N4").Value = 0.089589
$ws.Range("O4").Value = 0.005151837424850708
$ws.Range("P4").Value = 0.007275313219982055
$ws.Range("Q4").Value = 0.02194893668966667
$ws.Range("R4").Value = 0.197540430207
$ws.Range("S4").Value = 0.005151837424850708
$ws.Range("T4").Value = 0.007275313219982055

# --- Row 5 (new): target cluster = MuSCs (shifted down from row 4) --------------
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Vip"
$ws.Range("C5").Value = "Vipr2"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7349876666666667
$ws.Range("H5").Value = 2.204963
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.075609
$ws.Range("N5").Value = 10.151218
$ws.Range("O5").Value = 0.8756224223992591
$ws.Range("P5").Value = 0.8243566789931777
$ws.Range("Q5").Value = 3.730510015822333
$ws.Range("R5").Value = 22.383060094934
$ws.Range("S5").Value = 0.8756224223992591
$ws.Range("T5").Value = 0.8243566789931777
